# Weekly price-sheet update: a new week's Albahaca price record is inserted
# at row 25 ("Primera" quality, $300 volume band, $1300-1500 price range,
# date 2023-02-09 / serial 44966), pushing all subsequent data rows (old
# rows 25-52) down by one. The row that falls off the bottom (old row 52)
# becomes the new row 53.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 25-52 down to 26-53, creating a blank row 25
# (Excel automatically extends the used range / dimension).
$ws.Rows("25:25").Insert()

# Populate the newly inserted row 25 with the new weekly record.
$ws.Range("A25").Value = 1
$ws.Range("B25").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C25").Value = 'Arica y Parinacota'
$ws.Range("D25").Value = 44966
$ws.Range("E25").Value = 15
$ws.Range("F25").Value = 100112052
$ws.Range("G25").Value = 'Albahaca'
$ws.Range("H25").Value = 'Sin especificar'
$ws.Range("I25").Value = 'Primera'
$ws.Range("J25").Value = 300
$ws.Range("K25").Value = 1300
$ws.Range("L25").Value = 1500
$ws.Range("M25").Value = 1400
$ws.Range("N25").Value = '$/paquete'
$ws.Range("O25").Value = 'Región de Arica y Parinacota'
$ws.Range("P25").Value = 1400
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = 'Hortaliza'
